# "Generate Report for Archive"
#
# The localization status report is being regenerated:
#   - the handoff/status text moves from "Ready for handoff" to
#     "In Translation" everywhere it appears (Overview!E2:F3 and the
#     "Status" column on each language sheet), and
#   - the now-shorter status text lets the Status column(s) shrink, so
#     those columns are narrowed to fit the new content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E ("zh-cn") and F ("de-de") hold the status ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Re-fit the two status columns to the new (shorter) text, then match the
# precise target width used by the regenerated report.
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-language sheets: column C is "Status" ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Sheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    $ws.Columns.Item(3).AutoFit() | Out-Null
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
